$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 and 29 swap coins (Kaspa <-> WrappedeETH) along with their new
# price / volume figures.
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.510.34"
$ws.Range("E28").Value = "  -12.30%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.152"
$ws.Range("E29").Value = "  -11.83%  "

# Updated price (column D) and volume(1h) (column E) figures for the rest
# of the coin rows.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.448.90"
$ws.Range("E2").Value = "  -9.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.396.13"
$ws.Range("E3").Value = "  -12.12%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "461.35"
$ws.Range("E5").Value = "  -8.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.51"
$ws.Range("E6").Value = "  -7.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.486"
$ws.Range("E8").Value = "  -8.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.406.60"
$ws.Range("E9").Value = "  -12.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0942"
$ws.Range("E10").Value = "  -9.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("E11").Value = "  -12.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.312"
$ws.Range("E12").Value = "  -10.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.835.01"
$ws.Range("E14").Value = "  -11.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "53.468.44"
$ws.Range("E15").Value = "  -9.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000131"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.52"
$ws.Range("E17").Value = "  -9.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.428.18"
$ws.Range("E18").Value = "  -11.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.15"
$ws.Range("E19").Value = "  -13.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "305.97"
$ws.Range("E20").Value = "  -11.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.34"
$ws.Range("E21").Value = "  -15.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.68"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.32"
$ws.Range("E24").Value = "  -15.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.71"
$ws.Range("E25").Value = "  -11.87%  "
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.384"
$ws.Range("E27").Value = "  -10.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.04"
$ws.Range("E30").Value = "  -6.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0713"
$ws.Range("E32").Value = "  -14.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "145.76"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.58"
$ws.Range("E34").Value = "  -8.41%  "
$ws.Range("E35").Value = "  -11.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -8.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.49"
$ws.Range("E37").Value = "  -16.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.05"
$ws.Range("E38").Value = "  -6.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.795"
$ws.Range("E39").Value = "  -16.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.991"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "32.59"
$ws.Range("E41").Value = "  -9.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.588"
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0520"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.23"
$ws.Range("E44").Value = "  -8.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("E46").Value = "  -12.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.928.81"
$ws.Range("E47").Value = "  -11.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0867"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("E49").Value = "  -4.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.15"
$ws.Range("E50").Value = "  -12.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.33"
$ws.Range("E51").Value = "  -14.28%  "
